# Auto-generated Excel COM-interop script applying the Shiva_Profits value update.
# Updates specific leve-profit calculation cells (columns H-N) across all 8 job sheets,
# reflecting refreshed market-board price data from the scheduled scraper run.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H17").Value = 1222.0564
$ws.Range("J17").Value = 1210.9429
$ws.Range("L17").Value = 3632.8287
$ws.Range("N17").Value = -3968.8287
$ws.Range("H87").Value = 49998.332
$ws.Range("J87").Value = 49998.332
$ws.Range("L87").Value = 49998.332
$ws.Range("N87").Value = -52494.332
$ws.Range("H90").Value = 49998.332
$ws.Range("J90").Value = 49998.332
$ws.Range("L90").Value = 149994.996
$ws.Range("N90").Value = -162474.996
$ws.Range("H92").Value = 17544132
$ws.Range("I92").Value = 25641242
$ws.Range("J92").Value = 392.33334
$ws.Range("K92").Value = 25641242
$ws.Range("L92").Value = 392.33334
$ws.Range("M92").Value = -25639994
$ws.Range("N92").Value = -2888.33334
$ws.Range("H127").Value = 790.25
$ws.Range("I127").Value = 546
$ws.Range("K127").Value = 1638
$ws.Range("M127").Value = 3322
$ws.Range("H138").Value = 22223986
$ws.Range("I138").Value = 28572430
$ws.Range("J138").Value = 4433
$ws.Range("K138").Value = 85717290
$ws.Range("L138").Value = 13299
$ws.Range("M138").Value = -85712150
$ws.Range("N138").Value = -23579
$ws.Range("H141").Value = 2531.0312
$ws.Range("I141").Value = 2764.5088
$ws.Range("K141").Value = 8293.526400000001
$ws.Range("M141").Value = -3113.526400000001

$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H32").Value = 2704.91
$ws.Range("I32").Value = 2701.9292
$ws.Range("K32").Value = 2701.9292
$ws.Range("M32").Value = -2414.9292
$ws.Range("H102").Value = 4029.7144
$ws.Range("I102").Value = 2252.3
$ws.Range("K102").Value = 2252.3
$ws.Range("M102").Value = -630.3000000000002

$ws = $wb.Worksheets.Item("BSM")

$ws.Range("H107").Value = 3008.2058
$ws.Range("I107").Value = 1907.1786
$ws.Range("J107").Value = 8146.3335
$ws.Range("K107").Value = 1907.1786
$ws.Range("L107").Value = 8146.3335
$ws.Range("M107").Value = 12.82140000000004
$ws.Range("N107").Value = -11986.3335
$ws.Range("H137").Value = 142499.33
$ws.Range("J137").Value = 142499.33
$ws.Range("L137").Value = 142499.33
$ws.Range("N137").Value = -152699.33

$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H17").Value = 4999
$ws.Range("J17").Value = 4999
$ws.Range("L17").Value = 4999
$ws.Range("N17").Value = -5347
$ws.Range("H31").Value = 2053.423
$ws.Range("I31").Value = 1683.8096
$ws.Range("J31").Value = 3605.8
$ws.Range("K31").Value = 1683.8096
$ws.Range("L31").Value = 3605.8
$ws.Range("M31").Value = -1388.8096
$ws.Range("N31").Value = -4195.8
$ws.Range("H34").Value = 2053.423
$ws.Range("I34").Value = 1683.8096
$ws.Range("J34").Value = 3605.8
$ws.Range("K34").Value = 1683.8096
$ws.Range("L34").Value = 3605.8
$ws.Range("M34").Value = -1481.8096
$ws.Range("N34").Value = -4009.8
$ws.Range("H41").Value = 7992
$ws.Range("I41").Value = 3320
$ws.Range("K41").Value = 3320
$ws.Range("M41").Value = -2892
$ws.Range("H58").Value = 1648.1305
$ws.Range("J58").Value = 1895.4
$ws.Range("L58").Value = 1895.4
$ws.Range("N58").Value = -2301.4
$ws.Range("H105").Value = 3061.4614
$ws.Range("I105").Value = 1402.4546
$ws.Range("J105").Value = 12186
$ws.Range("K105").Value = 1402.4546
$ws.Range("L105").Value = 12186
$ws.Range("M105").Value = 344.5454
$ws.Range("N105").Value = -15680
$ws.Range("H122").Value = 3045.6155
$ws.Range("I122").Value = 2952.611
$ws.Range("K122").Value = 8857.832999999999
$ws.Range("M122").Value = -6407.832999999999
$ws.Range("H134").Value = 4001.3438
$ws.Range("I134").Value = 4044.8462
$ws.Range("J134").Value = 3812.8333
$ws.Range("K134").Value = 12134.5386
$ws.Range("L134").Value = 11438.4999
$ws.Range("M134").Value = -9599.5386
$ws.Range("N134").Value = -16508.4999
$ws.Range("H135").Value = 140699.2
$ws.Range("J135").Value = 140699.2
$ws.Range("L135").Value = 140699.2
$ws.Range("N135").Value = -150839.2
$ws.Range("H136").Value = 1648.1305
$ws.Range("J136").Value = 1895.4
$ws.Range("L136").Value = 5686.200000000001
$ws.Range("N136").Value = -10786.2

$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H13").Value = 2176
$ws.Range("I13").Value = 114.28571
$ws.Range("J13").Value = 5784
$ws.Range("K13").Value = 342.85713
$ws.Range("L13").Value = 17352
$ws.Range("M13").Value = -174.85713
$ws.Range("N13").Value = -17688
$ws.Range("H17").Value = 2587.5
$ws.Range("I17").Value = 2450
$ws.Range("J17").Value = 3000
$ws.Range("K17").Value = 7350
$ws.Range("L17").Value = 9000
$ws.Range("M17").Value = -7181
$ws.Range("N17").Value = -9338
$ws.Range("H19").Value = 1074.1666
$ws.Range("I19").Value = 579
$ws.Range("J19").Value = 2064.5
$ws.Range("K19").Value = 1737
$ws.Range("L19").Value = 6193.5
$ws.Range("M19").Value = -1563
$ws.Range("N19").Value = -6541.5
$ws.Range("H37").Value = 79401
$ws.Range("J37").Value = 79401
$ws.Range("L37").Value = 238203
$ws.Range("N37").Value = -238427
$ws.Range("H98").Value = 416
$ws.Range("I98").Value = 394.14285
$ws.Range("J98").Value = 454.25
$ws.Range("K98").Value = 1182.42855
$ws.Range("L98").Value = 1362.75
$ws.Range("M98").Value = 315.5714499999999
$ws.Range("N98").Value = -4358.75
$ws.Range("H107").Value = 399.69232
$ws.Range("J107").Value = 463.22223
$ws.Range("L107").Value = 1389.66669
$ws.Range("N107").Value = -5229.66669
$ws.Range("H134").Value = 1394.129
$ws.Range("I134").Value = 1394.129
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 4182.387
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = 887.6130000000003
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H70").Value = 9907.611000000001
$ws.Range("I70").Value = 10961.866
$ws.Range("J70").Value = 4636.3335
$ws.Range("K70").Value = 10961.866
$ws.Range("L70").Value = 4636.3335
$ws.Range("M70").Value = -10691.866
$ws.Range("N70").Value = -5176.3335
$ws.Range("H73").Value = 9907.611000000001
$ws.Range("I73").Value = 10961.866
$ws.Range("J73").Value = 4636.3335
$ws.Range("K73").Value = 10961.866
$ws.Range("L73").Value = 4636.3335
$ws.Range("M73").Value = -10025.866
$ws.Range("N73").Value = -6508.3335
$ws.Range("H96").Value = 100130.5
$ws.Range("J96").Value = 100130.5
$ws.Range("L96").Value = 100130.5
$ws.Range("N96").Value = -105622.5
$ws.Range("H102").Value = 7953
$ws.Range("I102").Value = 5688.778
$ws.Range("K102").Value = 5688.778
$ws.Range("M102").Value = -4066.778
$ws.Range("H122").Value = 1403.7959
$ws.Range("I122").Value = 1091.1951
$ws.Range("K122").Value = 3273.5853
$ws.Range("M122").Value = -823.5852999999997

$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H22").Value = 909.625
$ws.Range("I22").Value = 896.8
$ws.Range("J22").Value = 931
$ws.Range("K22").Value = 896.8
$ws.Range("L22").Value = 931
$ws.Range("M22").Value = -601.8
$ws.Range("N22").Value = -1521
$ws.Range("H27").Value = 909.625
$ws.Range("I27").Value = 896.8
$ws.Range("J27").Value = 931
$ws.Range("K27").Value = 896.8
$ws.Range("L27").Value = 931
$ws.Range("M27").Value = -789.8
$ws.Range("N27").Value = -1145
$ws.Range("H93").Value = 9091670
$ws.Range("I93").Value = 10000779
$ws.Range("K93").Value = 10000779
$ws.Range("M93").Value = -9999531
$ws.Range("H100").Value = 40002610
$ws.Range("I100").Value = 90911304
$ws.Range("K100").Value = 90911304
$ws.Range("M100").Value = -90910763
$ws.Range("H122").Value = 4176.5107
$ws.Range("I122").Value = 4171.6523
$ws.Range("K122").Value = 12514.9569
$ws.Range("M122").Value = -10064.9569
$ws.Range("H132").Value = 30759.91
$ws.Range("I132").Value = 31856.592
$ws.Range("K132").Value = 95569.776
$ws.Range("M132").Value = -93039.776
$ws.Range("H136").Value = 6391.926
$ws.Range("I136").Value = 6399.1055
$ws.Range("J136").Value = 6374.875
$ws.Range("K136").Value = 19197.3165
$ws.Range("L136").Value = 19124.625
$ws.Range("M136").Value = -16647.3165
$ws.Range("N136").Value = -24224.625

$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H122").Value = 3579.8545
$ws.Range("I122").Value = 1799.9375
$ws.Range("J122").Value = 15785
$ws.Range("K122").Value = 5399.8125
$ws.Range("L122").Value = 47355
$ws.Range("M122").Value = -2949.8125
$ws.Range("N122").Value = -52255
$ws.Range("H132").Value = 2996.923
$ws.Range("I132").Value = 2405.024
$ws.Range("K132").Value = 7215.072
$ws.Range("M132").Value = -4685.072
$ws.Range("H139").Value = 69949
$ws.Range("J139").Value = 69949
$ws.Range("L139").Value = 69949
$ws.Range("N139").Value = -80229
